$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 7.83870406380029
$ws.Range("D2").Value = 8.893515367309677
$ws.Range("E2").Value = 13.67117341556229
$ws.Range("F2").Value = 34.53481326657113
$ws.Range("G2").Value = 3.667112050585561
$ws.Range("J2").Value = 10.23893682429091
$ws.Range("K2").Value = 11.0324791635777
$ws.Range("M2").Value = 15.45962279084162
$ws.Range("O2").Value = 26.17731385551513

$ws.Range("B3").Value = 7.76919519555176
$ws.Range("D3").Value = 8.870513336643558
$ws.Range("E3").Value = 13.6833507156601
$ws.Range("F3").Value = 34.62120944210999
$ws.Range("G3").Value = 3.669182108272986
$ws.Range("J3").Value = 10.26554340961156
$ws.Range("K3").Value = 10.6798326832802
$ws.Range("M3").Value = 15.32734605329021
$ws.Range("O3").Value = 26.27123239826054

$ws.Range("B4").Value = 7.728061420991882
$ws.Range("D4").Value = 8.857608266166796
$ws.Range("E4").Value = 13.6931147846978
$ws.Range("F4").Value = 34.68275544210056
$ws.Range("G4").Value = 3.670520768706877
$ws.Range("J4").Value = 10.28309531660535
$ws.Range("K4").Value = 10.45805721172665
$ws.Range("M4").Value = 15.24755896922338
$ws.Range("O4").Value = 26.33482668801035

$ws.Range("B5").Value = 7.711705845019362
$ws.Range("D5").Value = 8.852659228893371
$ws.Range("E5").Value = 13.69766918126728
$ws.Range("F5").Value = 34.70996740789044
$ws.Range("G5").Value = 3.671083345492721
$ws.Range("J5").Value = 10.29055380660896
$ws.Range("K5").Value = 10.36649090765778
$ws.Range("M5").Value = 15.21543278662956
$ws.Range("O5").Value = 26.36222901546017

$ws.Range("B6").Value = 7.709015101416443
$ws.Range("D6").Value = 8.851856259924409
$ws.Range("E6").Value = 13.69846019853332
$ws.Range("F6").Value = 34.71461449891133
$ws.Range("G6").Value = 3.671177793018511
$ws.Range("J6").Value = 10.2918107713986
$ws.Range("K6").Value = 10.35121843728223
$ws.Range("M6").Value = 15.21012247315819
$ws.Range("O6").Value = 26.36686886672091

$ws.Range("B7").Value = 7.727839173376633
$ws.Range("D7").Value = 8.857540262508405
$ws.Range("E7").Value = 13.69317387653407
$ws.Range("F7").Value = 34.68311380962412
$ws.Range("G7").Value = 3.670528286659484
$ws.Range("J7").Value = 10.28319466520414
$ws.Range("K7").Value = 10.45682696060068
$ws.Range("M7").Value = 15.24712409831204
$ws.Range("O7").Value = 26.33519022872217

$ws.Range("B8").Value = 7.814427942580493
$ws.Range("D8").Value = 8.885333608738563
$ws.Range("E8").Value = 13.67489769122509
$ws.Range("F8").Value = 34.56283564686009
$ws.Range("G8").Value = 3.667811801053773
$ws.Range("J8").Value = 10.24785872906435
$ws.Range("K8").Value = 10.91205013658836
$ws.Range("M8").Value = 15.41373354376206
$ws.Range("O8").Value = 26.20846456199538

$ws.Range("B9").Value = 7.995649853199052
$ws.Range("D9").Value = 8.949337693630181
$ws.Range("E9").Value = 13.65718634129965
$ws.Range("F9").Value = 34.39462919560734
$ws.Range("G9").Value = 3.663018982600109
$ws.Range("J9").Value = 10.18819340132176
$ws.Range("K9").Value = 11.75797394293683
$ws.Range("M9").Value = 15.75056498986559
$ws.Range("O9").Value = 26.00714357662969

$ws.Range("B10").Value = 8.13459649917764
$ws.Range("D10").Value = 9.001919277900175
$ws.Range("E10").Value = 13.65519218700706
$ws.Range("F10").Value = 34.31257670598893
$ws.Range("G10").Value = 3.659819911460328
$ws.Range("J10").Value = 10.15020642698078
$ws.Range("K10").Value = 12.34486230539351
$ws.Range("M10").Value = 16.00250143485905
$ws.Range("O10").Value = 25.88820252024718

$ws.Range("B11").Value = 8.198812470240417
$ws.Range("D11").Value = 9.026992259225084
$ws.Range("E11").Value = 13.65666725036266
$ws.Range("F11").Value = 34.28431553134739
$ws.Range("G11").Value = 3.658433802429006
$ws.Range("J11").Value = 10.13419094703959
$ws.Range("K11").Value = 12.60325394291352
$ws.Range("M11").Value = 16.11771770698495
$ws.Range("O11").Value = 25.84042472247657

$ws.Range("B12").Value = 8.223252808112653
$ws.Range("D12").Value = 9.036647621785637
$ws.Range("E12").Value = 13.65756718856999
$ws.Range("F12").Value = 34.27492027314914
$ws.Range("G12").Value = 3.657918809125714
$ws.Range("J12").Value = 10.1283078860968
$ws.Range("K12").Value = 12.69978796101947
$ws.Range("M12").Value = 16.16140386658203
$ws.Range("O12").Value = 25.82324583331905

$ws.Range("B13").Value = 8.217984015319232
$ws.Range("D13").Value = 9.034561099483016
$ws.Range("E13").Value = 13.65735820944821
$ws.Range("F13").Value = 34.27688555448655
$ws.Range("G13").Value = 3.65802928280228
$ws.Range("J13").Value = 10.12956683338172
$ws.Range("K13").Value = 12.67905719732851
$ws.Range("M13").Value = 16.15199326532808
$ws.Range("O13").Value = 25.82690493244106

$ws.Range("B14").Value = 8.200820848809537
$ws.Range("D14").Value = 9.027783419479091
$ws.Range("E14").Value = 13.65673445485802
$ws.Range("F14").Value = 34.28351637536056
$ws.Range("G14").Value = 3.658391235570034
$ws.Range("J14").Value = 10.13370330488442
$ws.Range("K14").Value = 12.61122259292452
$ws.Range("M14").Value = 16.12131084496996
$ws.Range("O14").Value = 25.83899307750348

$ws.Range("B15").Value = 8.190323301621056
$ws.Range("D15").Value = 9.02365267740452
$ws.Range("E15").Value = 13.65639680430175
$ws.Range("F15").Value = 34.28774818841591
$ws.Range("G15").Value = 3.658614229389456
$ws.Range("J15").Value = 10.1362606624314
$ws.Range("K15").Value = 12.56949873278783
$ws.Range("M15").Value = 16.10252337603888
$ws.Range("O15").Value = 25.84651648211063

$ws.Range("B16").Value = 8.130418257782166
$ws.Range("D16").Value = 9.00030347815135
$ws.Range("E16").Value = 13.65514362618892
$ws.Range("F16").Value = 34.31460633014359
$ws.Range("G16").Value = 3.65991188445757
$ws.Range("J16").Value = 10.15127850960055
$ws.Range("K16").Value = 12.32779609792502
$ws.Range("M16").Value = 15.99498135563606
$ws.Range("O16").Value = 25.8914525640655

$ws.Range("B17").Value = 8.09391118965552
$ws.Range("D17").Value = 8.986271350235215
$ws.Range("E17").Value = 13.6549841892435
$ws.Range("F17").Value = 34.33340708962015
$ws.Range("G17").Value = 3.660725633033131
$ws.Range("J17").Value = 10.1608152870522
$ws.Range("K17").Value = 12.17726132218469
$ws.Range("M17").Value = 15.92914148265765
$ws.Range("O17").Value = 25.92064289947328

$ws.Range("B18").Value = 8.073009523177735
$ws.Range("D18").Value = 8.978309317821994
$ws.Range("E18").Value = 13.65511679328071
$ws.Range("F18").Value = 34.34507387581444
$ws.Range("G18").Value = 3.661200192722053
$ws.Range("J18").Value = 10.16641967035041
$ws.Range("K18").Value = 12.08987356621766
$ws.Range("M18").Value = 15.89133133188589
$ws.Range("O18").Value = 25.93802793685797

$ws.Range("B19").Value = 8.065949772242542
$ws.Range("D19").Value = 8.975632361656004
$ws.Range("E19").Value = 13.65520024938658
$ws.Range("F19").Value = 34.34917046232413
$ws.Range("G19").Value = 3.661361990618241
$ws.Range("J19").Value = 10.16833767819128
$ws.Range("K19").Value = 12.06015000335666
$ws.Range("M19").Value = 15.87854058971736
$ws.Range("O19").Value = 25.94401638807025

$ws.Range("B20").Value = 8.097787627963745
$ws.Range("D20").Value = 8.987753862096467
$ws.Range("E20").Value = 13.65497795465467
$ws.Range("F20").Value = 34.33131740080418
$ws.Range("G20").Value = 3.660638334394288
$ws.Range("J20").Value = 10.15978775817476
$ws.Range("K20").Value = 12.19336982818066
$ws.Range("M20").Value = 15.9361443497352
$ws.Range("O20").Value = 25.91747387826112

$ws.Range("B21").Value = 8.20585891903171
$ws.Range("D21").Value = 9.029769866651357
$ws.Range("E21").Value = 13.65690841259067
$ws.Range("F21").Value = 34.2815332592369
$ws.Range("G21").Value = 3.658284653084513
$ws.Range("J21").Value = 10.13248339488104
$ws.Range("K21").Value = 12.6311834842987
$ws.Range("M21").Value = 16.13032173757779
$ws.Range("O21").Value = 25.83541767700594

$ws.Range("B22").Value = 8.277197415582771
$ws.Range("D22").Value = 9.058164344194907
$ws.Range("E22").Value = 13.66015912825919
$ws.Range("F22").Value = 34.25661349827761
$ws.Range("G22").Value = 3.656804045589646
$ws.Range("J22").Value = 10.11569711797697
$ws.Range("K22").Value = 12.90963491651909
$ws.Range("M22").Value = 16.25754141945871
$ws.Range("O22").Value = 25.78711540879822

$ws.Range("B23").Value = 8.23906518359375
$ws.Range("D23").Value = 9.042925891233251
$ws.Range("E23").Value = 13.6582426025236
$ws.Range("F23").Value = 34.26921582280124
$ws.Range("G23").Value = 3.657589014152829
$ws.Range("J23").Value = 10.12455948030338
$ws.Range("K23").Value = 12.76174699654819
$ws.Range("M23").Value = 16.18962349143048
$ws.Range("O23").Value = 25.81240680282772

$ws.Range("B24").Value = 8.096034818383515
$ws.Range("D24").Value = 8.987083290350164
$ws.Range("E24").Value = 13.65498007466366
$ws.Range("F24").Value = 34.33225947678131
$ws.Range("G24").Value = 3.660677781158078
$ws.Range("J24").Value = 10.16025192530041
$ws.Range("K24").Value = 12.18608979328261
$ws.Range("M24").Value = 15.932978220093
$ws.Range("O24").Value = 25.91890471427925

$ws.Range("B25").Value = 7.945519820309974
$ws.Range("D25").Value = 8.93102899146162
$ws.Range("E25").Value = 13.66003974150889
$ws.Range("F25").Value = 34.43285881890529
$ws.Range("G25").Value = 3.664258735200466
$ws.Range("J25").Value = 10.20330580934919
$ws.Range("K25").Value = 11.53480656691877
$ws.Range("M25").Value = 15.65853474622108
$ws.Range("O25").Value = 26.05653308694428
